$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1545.6
$ws.Range("I28").Value = 1384.2307
$ws.Range("K28").Value = 1384.2307
$ws.Range("M28").Value = -899.2307000000001
$ws.Range("H64").Value = 6345
$ws.Range("I64").Value = 4741.6665
$ws.Range("K64").Value = 4741.6665
$ws.Range("M64").Value = -4493.6665
$ws.Range("H67").Value = 6345
$ws.Range("I67").Value = 4741.6665
$ws.Range("K67").Value = 4741.6665
$ws.Range("M67").Value = -3883.6665
$ws.Range("H87").Value = 124709.5
$ws.Range("J87").Value = 124709.5
$ws.Range("L87").Value = 124709.5
$ws.Range("N87").Value = -127205.5
$ws.Range("H90").Value = 124709.5
$ws.Range("J90").Value = 124709.5
$ws.Range("L90").Value = 374128.5
$ws.Range("N90").Value = -386608.5
$ws.Range("H127").Value = 10000
$ws.Range("I127").Value = 10000
$ws.Range("J127").Value = 10000
$ws.Range("K127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("M127").Value = -25040
$ws.Range("N127").Value = -39920
$ws.Range("H137").Value = 1697.3334
$ws.Range("I137").Value = 1499
$ws.Range("K137").Value = 4497
$ws.Range("M137").Value = -1947
$ws.Range("H138").Value = 3315.6667
$ws.Range("J138").Value = 4750
$ws.Range("L138").Value = 14250
$ws.Range("N138").Value = -24530
$ws.Range("H140").Value = 40780
$ws.Range("J140").Value = 40780
$ws.Range("L140").Value = 40780
$ws.Range("N140").Value = -51140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 20661.666
$ws.Range("I28").Value = 20661.666
$ws.Range("K28").Value = 20661.666
$ws.Range("M28").Value = -20469.666
$ws.Range("H61").Value = 5035.4443
$ws.Range("I61").Value = 3808
$ws.Range("J61").Value = 6964.2856
$ws.Range("K61").Value = 3808
$ws.Range("L61").Value = 6964.2856
$ws.Range("M61").Value = -3596
$ws.Range("N61").Value = -7388.2856
$ws.Range("H74").Value = 1820
$ws.Range("I74").Value = 1902.7333
$ws.Range("J74").Value = 1199.5
$ws.Range("K74").Value = 1902.7333
$ws.Range("L74").Value = 1199.5
$ws.Range("M74").Value = -1028.7333
$ws.Range("N74").Value = -2947.5
$ws.Range("H77").Value = 1820
$ws.Range("I77").Value = 1902.7333
$ws.Range("J77").Value = 1199.5
$ws.Range("K77").Value = 9513.666500000001
$ws.Range("L77").Value = 5997.5
$ws.Range("M77").Value = -5145.666500000001
$ws.Range("N77").Value = -14733.5
$ws.Range("H99").Value = 20661.666
$ws.Range("I99").Value = 20661.666
$ws.Range("K99").Value = 20661.666
$ws.Range("M99").Value = -17666.666
$ws.Range("H110").Value = 1588.6
$ws.Range("I110").Value = 1588.6
$ws.Range("K110").Value = 1588.6
$ws.Range("M110").Value = 456.4000000000001
$ws.Range("H122").Value = 3422.6667
$ws.Range("I122").Value = 3422.6667
$ws.Range("K122").Value = 10268.0001
$ws.Range("M122").Value = -7818.000100000001
$ws.Range("H132").Value = 2377.6365
$ws.Range("I132").Value = 2412.3684
$ws.Range("J132").Value = 2157.6667
$ws.Range("K132").Value = 7237.1052
$ws.Range("L132").Value = 6473.000100000001
$ws.Range("M132").Value = -4707.1052
$ws.Range("N132").Value = -11533.0001
$ws.Range("H136").Value = 5035.4443
$ws.Range("I136").Value = 3808
$ws.Range("J136").Value = 6964.2856
$ws.Range("K136").Value = 11424
$ws.Range("L136").Value = 20892.8568
$ws.Range("M136").Value = -8874
$ws.Range("N136").Value = -25992.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 9027
$ws.Range("I75").Value = 9027
$ws.Range("K75").Value = 9027
$ws.Range("M75").Value = -8091
$ws.Range("H78").Value = 9027
$ws.Range("I78").Value = 9027
$ws.Range("K78").Value = 27081
$ws.Range("M78").Value = -22401
$ws.Range("H134").Value = 3999.4443
$ws.Range("I134").Value = 4000
$ws.Range("J134").Value = 3997.5
$ws.Range("K134").Value = 12000
$ws.Range("L134").Value = 11992.5
$ws.Range("M134").Value = -9465
$ws.Range("N134").Value = -17062.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4414.418
$ws.Range("I31").Value = 2080.5715
$ws.Range("K31").Value = 2080.5715
$ws.Range("M31").Value = -1785.5715
$ws.Range("H34").Value = 4414.418
$ws.Range("I34").Value = 2080.5715
$ws.Range("K34").Value = 2080.5715
$ws.Range("M34").Value = -1878.5715
$ws.Range("H58").Value = 1974.8462
$ws.Range("J58").Value = 4147.75
$ws.Range("L58").Value = 4147.75
$ws.Range("N58").Value = -4553.75
$ws.Range("H62").Value = 1950
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 1950
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("H86").Value = 7238.364
$ws.Range("I86").Value = 6951.5
$ws.Range("K86").Value = 6951.5
$ws.Range("M86").Value = -5828.5
$ws.Range("H89").Value = 7238.364
$ws.Range("I89").Value = 6951.5
$ws.Range("K89").Value = 34757.5
$ws.Range("M89").Value = -29141.5
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = $null
$ws.Range("H134").Value = 404.51852
$ws.Range("I134").Value = 381.6154
$ws.Range("K134").Value = 1144.8462
$ws.Range("M134").Value = 1390.1538
$ws.Range("H136").Value = 1974.8462
$ws.Range("J136").Value = 4147.75
$ws.Range("L136").Value = 12443.25
$ws.Range("N136").Value = -17543.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1107.1
$ws.Range("J34").Value = 1733.1666
$ws.Range("L34").Value = 5199.4998
$ws.Range("N34").Value = -5367.4998
$ws.Range("H75").Value = 1850
$ws.Range("J75").Value = 1850
$ws.Range("L75").Value = 5550
$ws.Range("N75").Value = -7546
$ws.Range("H78").Value = 1850
$ws.Range("J78").Value = 1850
$ws.Range("L78").Value = 16650
$ws.Range("N78").Value = -26634
$ws.Range("H86").Value = 533.3333
$ws.Range("I86").Value = 200
$ws.Range("J86").Value = 700
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 2100
$ws.Range("M86").Value = 586
$ws.Range("N86").Value = -4472
$ws.Range("H87").Value = 5950
$ws.Range("I87").Value = 5950
$ws.Range("K87").Value = 17850
$ws.Range("M87").Value = -16602
$ws.Range("H89").Value = 533.3333
$ws.Range("I89").Value = 200
$ws.Range("J89").Value = 700
$ws.Range("K89").Value = 1800
$ws.Range("L89").Value = 6300
$ws.Range("M89").Value = 4128
$ws.Range("N89").Value = -18156
$ws.Range("H90").Value = 5950
$ws.Range("I90").Value = 5950
$ws.Range("K90").Value = 53550
$ws.Range("M90").Value = -47310
$ws.Range("H107").Value = 166.25
$ws.Range("I107").Value = 233.5
$ws.Range("J107").Value = 99
$ws.Range("K107").Value = 700.5
$ws.Range("L107").Value = 297
$ws.Range("M107").Value = 1219.5
$ws.Range("N107").Value = -4137
$ws.Range("H131").Value = 1742.9459
$ws.Range("I131").Value = 608.7692
$ws.Range("K131").Value = 1826.3076
$ws.Range("M131").Value = 3213.6924

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 6254.4165
$ws.Range("I14").Value = 3484
$ws.Range("J14").Value = 9024.833000000001
$ws.Range("K14").Value = 3484
$ws.Range("L14").Value = 9024.833000000001
$ws.Range("M14").Value = -3316
$ws.Range("N14").Value = -9360.833000000001
$ws.Range("H136").Value = 22500
$ws.Range("J136").Value = 22500
$ws.Range("L136").Value = 67500
$ws.Range("N136").Value = -72600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6207.7144
$ws.Range("I7").Value = 5439.8
$ws.Range("K7").Value = 5439.8
$ws.Range("M7").Value = -5327.8
$ws.Range("H64").Value = 26500
$ws.Range("J64").Value = 26500
$ws.Range("L64").Value = 26500
$ws.Range("N64").Value = -26950
$ws.Range("H67").Value = 26500
$ws.Range("J67").Value = 26500
$ws.Range("L67").Value = 26500
$ws.Range("N67").Value = -28060
$ws.Range("H68").Value = 8247.875
$ws.Range("I68").Value = 7330.3335
$ws.Range("J68").Value = 8798.4
$ws.Range("K68").Value = 7330.3335
$ws.Range("L68").Value = 8798.4
$ws.Range("M68").Value = -6581.3335
$ws.Range("N68").Value = -10296.4
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
$ws.Range("H71").Value = 8247.875
$ws.Range("I71").Value = 7330.3335
$ws.Range("J71").Value = 8798.4
$ws.Range("K71").Value = 36651.6675
$ws.Range("L71").Value = 43992
$ws.Range("M71").Value = -32907.6675
$ws.Range("N71").Value = -51480
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
$ws.Range("H126").Value = 6207.7144
$ws.Range("I126").Value = 5439.8
$ws.Range("K126").Value = 16319.4
$ws.Range("M126").Value = -13849.4
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3613
$ws.Range("I122").Value = 3613
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10839
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8389
$ws.Range("N122").Value = $null
